$wb = $excel.ActiveWorkbook

# ALC row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 2232680.2
$ws.Range("I92").Value = 919646.4399999999
$ws.Range("K92").Value = 919646.4399999999
$ws.Range("M92").Value = -918398.4399999999

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1430
$ws.Range("I100").Value = 1430
$ws.Range("K100").Value = 1430
$ws.Range("M100").Value = -889

# ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 947.5
$ws.Range("I111").Value = 947.5
$ws.Range("K111").Value = 2842.5
$ws.Range("M111").Value = 224.5

# ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1624981.6
$ws.Range("I125").Value = 3789312
$ws.Range("J125").Value = 1733.75
$ws.Range("K125").Value = 34103808
$ws.Range("L125").Value = 15603.75
$ws.Range("M125").Value = -34101348
$ws.Range("N125").Value = -20523.75

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1510.75
$ws.Range("I2").Value = 1426.5714
$ws.Range("K2").Value = 1426.5714
$ws.Range("M2").Value = -1313.5714

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 49692.715
$ws.Range("I45").Value = 60707.53
$ws.Range("K45").Value = 60707.53
$ws.Range("M45").Value = -60330.53

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 926411.9399999999
$ws.Range("I61").Value = 26566.488
$ws.Range("J61").Value = 2854652.2
$ws.Range("K61").Value = 26566.488
$ws.Range("L61").Value = 2854652.2
$ws.Range("M61").Value = -26354.488
$ws.Range("N61").Value = -2855076.2

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 339354.84
$ws.Range("I74").Value = 2233.2173
$ws.Range("K74").Value = 2233.2173
$ws.Range("M74").Value = -1359.2173

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 339354.84
$ws.Range("I77").Value = 2233.2173
$ws.Range("K77").Value = 11166.0865
$ws.Range("M77").Value = -6798.086499999999

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1510.75
$ws.Range("I116").Value = 1426.5714
$ws.Range("K116").Value = 1426.5714
$ws.Range("M116").Value = 867.4286

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 926411.9399999999
$ws.Range("I136").Value = 26566.488
$ws.Range("J136").Value = 2854652.2
$ws.Range("K136").Value = 79699.46400000001
$ws.Range("L136").Value = 8563956.600000001
$ws.Range("M136").Value = -77149.46400000001
$ws.Range("N136").Value = -8569056.600000001

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1510.75
$ws.Range("I3").Value = 1426.5714
$ws.Range("K3").Value = 1426.5714
$ws.Range("M3").Value = -1312.5714

# BSM row 61
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 19000
$ws.Range("J61").Value = 19000
$ws.Range("L61").Value = 19000
$ws.Range("N61").Value = -19626

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 30798.4
$ws.Range("I99").Value = 30798.4
$ws.Range("K99").Value = 30798.4
$ws.Range("M99").Value = -29300.4

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 7768.7095
$ws.Range("I105").Value = 6674.227
$ws.Range("J105").Value = 10444.111
$ws.Range("K105").Value = 6674.227
$ws.Range("L105").Value = 10444.111
$ws.Range("M105").Value = -4927.227
$ws.Range("N105").Value = -13938.111

# BSM row 135
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 88991.336
$ws.Range("J135").Value = 88991.336
$ws.Range("L135").Value = 88991.336
$ws.Range("N135").Value = -99131.336

# BSM row 138
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 153330.83
$ws.Range("J138").Value = 153330.83
$ws.Range("L138").Value = 153330.83
$ws.Range("N138").Value = -163610.83

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2428.6
$ws.Range("I58").Value = 1381
$ws.Range("K58").Value = 1381
$ws.Range("M58").Value = -1178

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2380.8948
$ws.Range("I107").Value = 2424.9285
$ws.Range("J107").Value = 2257.6
$ws.Range("K107").Value = 2424.9285
$ws.Range("L107").Value = 2257.6
$ws.Range("M107").Value = -504.9285
$ws.Range("N107").Value = -6097.6

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2813.0625
$ws.Range("I134").Value = 2546.2727
$ws.Range("J134").Value = 3400
$ws.Range("K134").Value = 7638.8181
$ws.Range("L134").Value = 10200
$ws.Range("M134").Value = -5103.8181
$ws.Range("N134").Value = -15270

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2428.6
$ws.Range("I136").Value = 1381
$ws.Range("K136").Value = 4143
$ws.Range("M136").Value = -1593

# CRP row 138
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 92527.30499999999
$ws.Range("J138").Value = 93571.25
$ws.Range("L138").Value = 93571.25
$ws.Range("N138").Value = -103851.25

# CRP row 140
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 426661.66
$ws.Range("J140").Value = 426661.66
$ws.Range("L140").Value = 426661.66
$ws.Range("N140").Value = -437021.66

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1710.909
$ws.Range("I5").Value = 1164.2
$ws.Range("K5").Value = 3492.6
$ws.Range("M5").Value = -3380.6

# CUL row 23
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1597.8
$ws.Range("I23").Value = 3594.5
$ws.Range("K23").Value = 10783.5
$ws.Range("M23").Value = -10548.5

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1710.909
$ws.Range("I135").Value = 1164.2
$ws.Range("K135").Value = 10477.8
$ws.Range("M135").Value = -7942.800000000001

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 143927.14
$ws.Range("I107").Value = 500247
$ws.Range("K107").Value = 500247
$ws.Range("M107").Value = -498327

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2005.0571
$ws.Range("I40").Value = 1969
$ws.Range("K40").Value = 1969
$ws.Range("M40").Value = -1833

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1725.5
$ws.Range("I82").Value = 1834
$ws.Range("J82").Value = 1400
$ws.Range("K82").Value = 1834
$ws.Range("L82").Value = 1400
$ws.Range("M82").Value = -1473
$ws.Range("N82").Value = -2122

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1725.5
$ws.Range("I85").Value = 1834
$ws.Range("J85").Value = 1400
$ws.Range("K85").Value = 1834
$ws.Range("L85").Value = 1400
$ws.Range("M85").Value = -586
$ws.Range("N85").Value = -3896

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3619.5
$ws.Range("I122").Value = 3197.1
$ws.Range("K122").Value = 9591.299999999999
$ws.Range("M122").Value = -7141.299999999999

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3215.389
$ws.Range("I132").Value = 2148.4285
$ws.Range("K132").Value = 6445.2855
$ws.Range("M132").Value = -3915.2855

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 58251.723
$ws.Range("I136").Value = 102815.3
$ws.Range("K136").Value = 308445.9
$ws.Range("M136").Value = -305895.9

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1682556
$ws.Range("I107").Value = 1510
$ws.Range("J107").Value = 5717066.5
$ws.Range("K107").Value = 4530
$ws.Range("L107").Value = 17151199.5
$ws.Range("M107").Value = -2610
$ws.Range("N107").Value = -17155039.5

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 666.41174
$ws.Range("I113").Value = 596.0909
$ws.Range("K113").Value = 1788.2727
$ws.Range("M113").Value = 381.7273

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1937.4546
$ws.Range("I122").Value = 1726.625
$ws.Range("J122").Value = 2499.6667
$ws.Range("K122").Value = 5179.875
$ws.Range("L122").Value = 7499.000100000001
$ws.Range("M122").Value = -2729.875
$ws.Range("N122").Value = -12399.0001

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 20834804
$ws.Range("I126").Value = 22728514
$ws.Range("K126").Value = 68185542
$ws.Range("M126").Value = -68183072

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 867.4375
$ws.Range("I136").Value = 757.1667
$ws.Range("K136").Value = 2271.5001
$ws.Range("M136").Value = 278.4998999999998
